# Applies the 2026-02-16 22:21 meteocat refresh: new DATA_EXTRACCIO timestamps
# plus the handful of observation values (humidity, wind gust, temp, etc.)
# that shifted between the 21:48-21:50 run and this 22:19-22:21 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-16 22:19:04"
$ws.Range("N2").Value = "0.5 °C 21:38 TU"
$ws.Range("O2").Value = "2.0 °C"
$ws.Range("E3").Value = "2026-02-16 22:19:07"
$ws.Range("N3").Value = "-5.1 °C 21:59 TU"
$ws.Range("O3").Value = "-1.1 °C"
$ws.Range("E4").Value = "2026-02-16 22:19:10"
$ws.Range("N4").Value = "6.9 °C 21:40 TU"
$ws.Range("O4").Value = "13.2 °C"
$ws.Range("E5").Value = "2026-02-16 22:19:12"
$ws.Range("N5").Value = "-4.9 °C 21:53 TU"
$ws.Range("O5").Value = "-0.9 °C"
$ws.Range("E6").Value = "2026-02-16 22:19:15"
$ws.Range("J6").Value = "1012.5 hPa"
$ws.Range("E7").Value = "2026-02-16 22:19:18"
$ws.Range("E8").Value = "2026-02-16 22:19:21"
$ws.Range("E9").Value = "2026-02-16 22:19:23"
$ws.Range("H9").Value = "'69%"
$ws.Range("O9").Value = "11.6 °C"
$ws.Range("E10").Value = "2026-02-16 22:19:26"
$ws.Range("O10").Value = "10.7 °C"
$ws.Range("E11").Value = "2026-02-16 22:19:29"
$ws.Range("O11").Value = "7.1 °C"
$ws.Range("E12").Value = "2026-02-16 22:19:31"
$ws.Range("H12").Value = "'77%"
$ws.Range("E13").Value = "2026-02-16 22:19:34"
$ws.Range("K13").Value = "12.6 MJ/m2"
$ws.Range("E14").Value = "2026-02-16 22:19:37"
$ws.Range("E15").Value = "2026-02-16 22:19:40"
$ws.Range("H15").Value = "'65%"
$ws.Range("O15").Value = "11.9 °C"
$ws.Range("E16").Value = "2026-02-16 22:19:42"
$ws.Range("N16").Value = "-4.1 °C 21:49 TU"
$ws.Range("E17").Value = "2026-02-16 22:19:45"
$ws.Range("E18").Value = "2026-02-16 22:19:48"
$ws.Range("J18").Value = "1012.8 hPa"
$ws.Range("O18").Value = "10.6 °C"
$ws.Range("E19").Value = "2026-02-16 22:19:51"
$ws.Range("O19").Value = "7.2 °C"
$ws.Range("E20").Value = "2026-02-16 22:19:53"
$ws.Range("L20").Value = "69.8 km/h - 336º 21:55 TU"
$ws.Range("N20").Value = "-2.9 °C 21:59 TU"
$ws.Range("E21").Value = "2026-02-16 22:19:56"
$ws.Range("H21").Value = "'68%"
$ws.Range("E22").Value = "2026-02-16 22:19:59"
$ws.Range("E23").Value = "2026-02-16 22:20:01"
$ws.Range("I23").Value = "16.2 mm"
$ws.Range("N23").Value = "-4.9 °C 21:56 TU"
$ws.Range("O23").Value = "-0.9 °C"
$ws.Range("E24").Value = "2026-02-16 22:20:04"
$ws.Range("E25").Value = "2026-02-16 22:20:07"
$ws.Range("H25").Value = "'86%"
$ws.Range("L25").Value = "33.5 km/h - 240º 21:42 TU"
$ws.Range("N25").Value = "-2.6 °C 21:57 TU"
$ws.Range("O25").Value = "0.5 °C"
$ws.Range("E26").Value = "2026-02-16 22:20:09"
$ws.Range("E27").Value = "2026-02-16 22:20:12"
$ws.Range("N27").Value = "-0.4 °C 21:44 TU"
$ws.Range("O27").Value = "1.1 °C"
$ws.Range("E28").Value = "2026-02-16 22:20:14"
$ws.Range("H28").Value = "'73%"
$ws.Range("E29").Value = "2026-02-16 22:20:17"
$ws.Range("L29").Value = "24.8 km/h - 336º 21:43 TU"
$ws.Range("O29").Value = "10.9 °C"
$ws.Range("E30").Value = "2026-02-16 22:20:20"
$ws.Range("E31").Value = "2026-02-16 22:20:22"
$ws.Range("J31").Value = "1011.7 hPa"
$ws.Range("O31").Value = "14.3 °C"
$ws.Range("E32").Value = "2026-02-16 22:20:25"
$ws.Range("L32").Value = "56.5 km/h - 277º 21:58 TU"
$ws.Range("E33").Value = "2026-02-16 22:20:28"
$ws.Range("E34").Value = "2026-02-16 22:20:31"
$ws.Range("N34").Value = "0.2 °C 21:55 TU"
$ws.Range("E35").Value = "2026-02-16 22:20:33"
$ws.Range("H35").Value = "'75%"
$ws.Range("I35").Value = "2.6 mm"
$ws.Range("O35").Value = "9.4 °C"
$ws.Range("E36").Value = "2026-02-16 22:20:36"
$ws.Range("E37").Value = "2026-02-16 22:20:39"
$ws.Range("E38").Value = "2026-02-16 22:20:41"
$ws.Range("O38").Value = "11.8 °C"
$ws.Range("E39").Value = "2026-02-16 22:20:44"
$ws.Range("N39").Value = "-4.9 °C 21:59 TU"
$ws.Range("O39").Value = "0.0 °C"
$ws.Range("E40").Value = "2026-02-16 22:20:47"
$ws.Range("H40").Value = "'84%"
$ws.Range("E41").Value = "2026-02-16 22:20:49"
$ws.Range("E42").Value = "2026-02-16 22:20:52"
$ws.Range("E43").Value = "2026-02-16 22:20:55"
$ws.Range("E44").Value = "2026-02-16 22:20:57"
$ws.Range("I44").Value = "13.8 mm"
$ws.Range("N44").Value = "-4.0 °C 21:59 TU"
$ws.Range("O44").Value = "-0.4 °C"
$ws.Range("E45").Value = "2026-02-16 22:21:00"
$ws.Range("H45").Value = "'95%"
$ws.Range("E46").Value = "2026-02-16 22:21:03"
